$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 68 of profit data (run date 2026-01-31)
$ws.Range("A68").NumberFormat = "@"
$ws.Range("A68").Value = "01/31/2026"
$ws.Range("A68").Style = "Normal"

$ws.Range("B68").Value = 10417.61
$ws.Range("C68").Value = 0.249843725829219
$ws.Range("D68").Value = 0.750156274170781
$ws.Range("E68").Value = -272.32
$ws.Range("F68").Value = -34.69
$ws.Range("G68").Value = -23026.92
$ws.Range("H68").Value = -74.66
$ws.Range("I68").Value = -650.08
$ws.Range("J68").Value = -19.98
$ws.Range("K68").Value = -23677
$ws.Range("L68").Value = -69.44
